# Trade #7 closed at 2026-02-18 00:09:46 - unknown UNKNOWN +0.000%
#
# Trade #35 (overall) / Trade #7 (MarketMaking strategy sheet) transitions
# from OPEN to CLOSED (early_exit), which updates the aggregate metrics on
# the Summary and Strategy Status sheets as well as the trade row itself on
# the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.82   # Current Capital
$summary.Range("B4").Value = 0.92      # Total P&L $
$summary.Range("B5").Value = 0.53      # Total P&L %
$summary.Range("B6").Value = 35        # Total Trades
$summary.Range("B7").Value = 19        # Winning Trades
$summary.Range("B9").Value = 54.29     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.82      # Capital
$status.Range("D6").Value = 6          # Trades
$status.Range("E6").Value = 0.01       # P&L $
$status.Range("F6").Value = -0.18      # P&L %
$status.Range("G6").Value = 50         # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet - Trade #35 (row 36)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G36").Value = 0.59           # Exit Price
$allTrades.Range("H36").Value = "CLOSED"       # Status
$allTrades.Range("I36").Value = 25.375         # P&L %
$allTrades.Range("J36").Value = 0.12           # P&L $
$allTrades.Range("K36").Value = 99.82          # Capital After
$allTrades.Range("L36").Value = "early_exit"   # Exit Reason
$allTrades.Range("M36").Value = 0.17           # Duration (min)

# ---------------------------------------------------------------------
# MarketMaking sheet - same trade, local row 7
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("G7").Value = 0.59         # Exit Price
$marketMaking.Range("H7").Value = "CLOSED"     # Status
$marketMaking.Range("I7").Value = 25.375       # P&L %
$marketMaking.Range("J7").Value = 0.12         # P&L $
$marketMaking.Range("K7").Value = 99.82        # Capital After
$marketMaking.Range("P7").Value = "early_exit" # Exit Reason
$marketMaking.Range("Q7").Value = 0.17         # Duration (min)
